$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 2613.1538
$ws.Range("I107").Value = 1942.8
$ws.Range("J107").Value = 4847.6665
$ws.Range("K107").Value = 1942.8
$ws.Range("L107").Value = 4847.6665
$ws.Range("M107").Value = -22.79999999999995
$ws.Range("N107").Value = -8687.666499999999

$ws.Range("H132").Value = 589436.75
$ws.Range("I132").Value = 1172.8572
$ws.Range("J132").Value = 3334668.2
$ws.Range("K132").Value = 3518.5716
$ws.Range("L132").Value = 10004004.6
$ws.Range("M132").Value = -988.5715999999998
$ws.Range("N132").Value = -10009064.6

$ws.Range("H135").Value = 10512
$ws.Range("I135").Value = 584.625
$ws.Range("K135").Value = 5261.625
$ws.Range("M135").Value = -2726.625

$ws.Range("H137").Value = 5979
$ws.Range("I137").Value = 4916.1665
$ws.Range("J137").Value = 6890
$ws.Range("K137").Value = 14748.4995
$ws.Range("L137").Value = 20670
$ws.Range("M137").Value = -12198.4995
$ws.Range("N137").Value = -25770

$ws.Range("H138").Value = 4085.9348
$ws.Range("I138").Value = 1885.0416
$ws.Range("J138").Value = 6486.909
$ws.Range("K138").Value = 5655.1248
$ws.Range("L138").Value = 19460.727
$ws.Range("M138").Value = -515.1247999999996
$ws.Range("N138").Value = -29740.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1037.28
$ws.Range("I2").Value = 1037.28
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1037.28
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -924.28
$ws.Range("N2").ClearContents()

$ws.Range("H32").Value = 5427.7812
$ws.Range("I32").Value = 2793.2788
$ws.Range("K32").Value = 2793.2788
$ws.Range("M32").Value = -2506.2788

$ws.Range("H61").Value = 3001.2083
$ws.Range("I61").Value = 2334.4375
$ws.Range("J61").Value = 4334.75
$ws.Range("K61").Value = 2334.4375
$ws.Range("L61").Value = 4334.75
$ws.Range("M61").Value = -2122.4375
$ws.Range("N61").Value = -4758.75

$ws.Range("H74").Value = 1250
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 1250
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 1250
$ws.Range("N74").Value = -2998
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 1250
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 1250
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 6250
$ws.Range("N77").Value = -14986
$ws.Range("M77").ClearContents()

$ws.Range("H102").Value = 3256.8333
$ws.Range("I102").Value = 1090.3478
$ws.Range("J102").Value = 7089.846
$ws.Range("K102").Value = 1090.3478
$ws.Range("L102").Value = 7089.846
$ws.Range("M102").Value = 531.6522
$ws.Range("N102").Value = -10333.846

$ws.Range("H116").Value = 1037.28
$ws.Range("I116").Value = 1037.28
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1037.28
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1256.72
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 1965.4231
$ws.Range("I132").Value = 1886.4166
$ws.Range("J132").Value = 2913.5
$ws.Range("K132").Value = 5659.2498
$ws.Range("L132").Value = 8740.5
$ws.Range("M132").Value = -3129.2498
$ws.Range("N132").Value = -13800.5

$ws.Range("H136").Value = 3001.2083
$ws.Range("I136").Value = 2334.4375
$ws.Range("J136").Value = 4334.75
$ws.Range("K136").Value = 7003.3125
$ws.Range("L136").Value = 13004.25
$ws.Range("M136").Value = -4453.3125
$ws.Range("N136").Value = -18104.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1037.28
$ws.Range("I3").Value = 1037.28
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1037.28
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -923.28
$ws.Range("N3").ClearContents()

$ws.Range("H134").Value = 3689.6875
$ws.Range("I134").Value = 3754.4614
$ws.Range("J134").Value = 3409
$ws.Range("K134").Value = 11263.3842
$ws.Range("L134").Value = 10227
$ws.Range("M134").Value = -8728.3842
$ws.Range("N134").Value = -15297

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 216.66667
$ws.Range("I7").Value = 216.66667
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 216.66667
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -103.66667
$ws.Range("N7").ClearContents()

$ws.Range("H31").Value = 7990.757
$ws.Range("I31").Value = 6613.3335
$ws.Range("J31").Value = 8433.5
$ws.Range("K31").Value = 6613.3335
$ws.Range("L31").Value = 8433.5
$ws.Range("M31").Value = -6318.3335
$ws.Range("N31").Value = -9023.5

$ws.Range("H34").Value = 7990.757
$ws.Range("I34").Value = 6613.3335
$ws.Range("J34").Value = 8433.5
$ws.Range("K34").Value = 6613.3335
$ws.Range("L34").Value = 8433.5
$ws.Range("M34").Value = -6411.3335
$ws.Range("N34").Value = -8837.5

$ws.Range("H132").Value = 1116.2273
$ws.Range("I132").Value = 951.97437
$ws.Range("K132").Value = 2855.92311
$ws.Range("M132").Value = -325.9231100000002

$ws.Range("H134").Value = 2758.1538
$ws.Range("I134").Value = 1906.4098
$ws.Range("J134").Value = 15747.25
$ws.Range("K134").Value = 5719.2294
$ws.Range("L134").Value = 47241.75
$ws.Range("M134").Value = -3184.2294
$ws.Range("N134").Value = -52311.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3520172.5
$ws.Range("I131").Value = 556005.5
$ws.Range("J131").Value = 3789642.2
$ws.Range("K131").Value = 1668016.5
$ws.Range("L131").Value = 11368926.6
$ws.Range("M131").Value = -1662976.5
$ws.Range("N131").Value = -11379006.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1479.6
$ws.Range("I113").Value = 1466
$ws.Range("K113").Value = 1466
$ws.Range("M113").Value = 704

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1964.4
$ws.Range("I22").Value = 1378.2858
$ws.Range("J22").Value = 3332
$ws.Range("K22").Value = 1378.2858
$ws.Range("L22").Value = 3332
$ws.Range("M22").Value = -1083.2858
$ws.Range("N22").Value = -3922

$ws.Range("H27").Value = 1964.4
$ws.Range("I27").Value = 1378.2858
$ws.Range("J27").Value = 3332
$ws.Range("K27").Value = 1378.2858
$ws.Range("L27").Value = 3332
$ws.Range("M27").Value = -1271.2858
$ws.Range("N27").Value = -3546

$ws.Range("H46").Value = 1295.7142
$ws.Range("I46").Value = 1184.75
$ws.Range("J46").Value = 1443.6666
$ws.Range("K46").Value = 1184.75
$ws.Range("L46").Value = 1443.6666
$ws.Range("M46").Value = -996.75
$ws.Range("N46").Value = -1819.6666

$ws.Range("H55").Value = 388
$ws.Range("I55").Value = 354.45456
$ws.Range("J55").Value = 461.8
$ws.Range("K55").Value = 354.45456
$ws.Range("L55").Value = 461.8
$ws.Range("M55").Value = -181.45456
$ws.Range("N55").Value = -807.8

$ws.Range("H61").Value = 2524.2144
$ws.Range("I61").Value = 2575.111
$ws.Range("J61").Value = 2432.6
$ws.Range("K61").Value = 2575.111
$ws.Range("L61").Value = 2432.6
$ws.Range("M61").Value = -2373.111
$ws.Range("N61").Value = -2836.6

$ws.Range("H113").Value = 2524.2144
$ws.Range("I113").Value = 2575.111
$ws.Range("J113").Value = 2432.6
$ws.Range("K113").Value = 2575.111
$ws.Range("L113").Value = 2432.6
$ws.Range("M113").Value = -405.1109999999999
$ws.Range("N113").Value = -6772.6

$ws.Range("H122").Value = 3339.2727
$ws.Range("I122").Value = 2822.5
$ws.Range("J122").Value = 3959.4
$ws.Range("K122").Value = 8467.5
$ws.Range("L122").Value = 11878.2
$ws.Range("M122").Value = -6017.5
$ws.Range("N122").Value = -16778.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3392.5715
$ws.Range("I126").Value = 2884.3845
$ws.Range("J126").Value = 9999
$ws.Range("K126").Value = 8653.1535
$ws.Range("L126").Value = 29997
$ws.Range("M126").Value = -6183.1535
$ws.Range("N126").Value = -34937

$ws.Range("H132").Value = 3711.9524
$ws.Range("I132").Value = 3230.257
$ws.Range("J132").Value = 6120.4287
$ws.Range("K132").Value = 9690.771000000001
$ws.Range("L132").Value = 18361.2861
$ws.Range("M132").Value = -7160.771000000001
$ws.Range("N132").Value = -23421.2861

$ws.Range("H136").Value = 610.3214
$ws.Range("I136").Value = 606.2273
$ws.Range("J136").Value = 625.3333
$ws.Range("K136").Value = 1818.6819
$ws.Range("L136").Value = 1875.9999
$ws.Range("M136").Value = 731.3181
$ws.Range("N136").Value = -6975.9999
